$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for column AD (index 30) and recomputed totals in column AG (index 33)
$updates = @(
    @{ Row = 2; AD = 10399.36; AG = 243940.58 },
    @{ Row = 3; AD = 5124.9;   AG = 124521.81 },
    @{ Row = 4; AD = 3131;     AG = 84636.89999999999 },
    @{ Row = 5; AD = 1877.5;   AG = 69839.28999999999 },
    @{ Row = 6; AD = 20532.76; AG = 522938.58 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 30).Value = $u.AD
    $ws.Cells.Item($u.Row, 33).Value = $u.AG
}
